$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# A new handoff xliff batch was generated for the four files that were
# "Ready for handoff" (31952d01, 42a3fa32, 6a77176b, 911ede20): their
# Priority flips from "low" to "ht" and the handoff timestamps advance.

# Overview sheet: "Latest HO Xliff Generate Date" for that batch advances.
$overview.Range("G4").Value = "2016-08-24 16:32:10"
$overview.Range("G5").Value = "2016-08-24 16:32:10"
$overview.Range("G6").Value = "2016-08-24 16:32:10"
$overview.Range("G7").Value = "2016-08-24 16:32:10"

# zh-cn sheet: Priority low -> ht, Latest Handoff Datetime advances.
$zhcn.Range("E4").Value = "ht"
$zhcn.Range("E5").Value = "ht"
$zhcn.Range("E6").Value = "ht"
$zhcn.Range("E7").Value = "ht"

$zhcn.Range("H4").Value = "2016-08-24 16:31:58"
$zhcn.Range("H5").Value = "2016-08-24 16:31:58"
$zhcn.Range("H6").Value = "2016-08-24 16:31:58"
$zhcn.Range("H7").Value = "2016-08-24 16:31:58"

# de-de sheet: Priority low -> ht. Its "Latest Handoff Datetime" cells
# happened to share the same text as Overview's date ("2016-08-24
# 16:31:42"), so they move in lockstep to the new value too.
$dede.Range("E4").Value = "ht"
$dede.Range("E5").Value = "ht"
$dede.Range("E6").Value = "ht"
$dede.Range("E7").Value = "ht"

$dede.Range("H4").Value = "2016-08-24 16:32:10"
$dede.Range("H5").Value = "2016-08-24 16:32:10"
$dede.Range("H6").Value = "2016-08-24 16:32:10"
$dede.Range("H7").Value = "2016-08-24 16:32:10"
